$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting so values
# like "1.002" / "12.20" are not silently coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.423.91'
$ws.Range("E2").Value = '  -2.90%  '
$ws.Range("D3").Value = '1.750.36'
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '321.91'
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -4.97%  '
$ws.Range("D8").Value = '0.3603'
$ws.Range("E8").Value = '  -2.68%  '
$ws.Range("D9").Value = '0.07513'
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("D10").Value = '42.18'
$ws.Range("E10").Value = '  -5.78%  '
$ws.Range("D11").Value = '1.102'
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("E13").Value = '  -6.34%  '
$ws.Range("D14").Value = '6.031'
$ws.Range("E14").Value = '  -4.23%  '
$ws.Range("D15").Value = '7.219'
$ws.Range("E15").Value = '  -4.61%  '
$ws.Range("D16").Value = '1.749.24'
$ws.Range("E16").Value = '  -5.26%  '
$ws.Range("D17").Value = '93.28'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '0.00001069'
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("D19").Value = '0.06355'
$ws.Range("E19").Value = '  -3.16%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -2.96%  '
$ws.Range("D22").Value = '5.885'
$ws.Range("E22").Value = '  -5.52%  '
$ws.Range("D23").Value = '27.480.39'
$ws.Range("E23").Value = '  -2.90%  '
$ws.Range("E24").Value = '  -4.16%  '
$ws.Range("D25").Value = '2.086'
$ws.Range("E25").Value = '  -2.86%  '
$ws.Range("D26").Value = '162.10'
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("D27").Value = '20.30'
$ws.Range("E27").Value = '  -2.29%  '
$ws.Range("D28").Value = '1.945.57'
$ws.Range("E28").Value = '  -4.52%  '
$ws.Range("D29").Value = '2.130'
$ws.Range("E29").Value = '  -8.07%  '
$ws.Range("D30").Value = '123.85'
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("E31").Value = '  -7.81%  '
$ws.Range("D32").Value = '3.644'
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").Value = '5.533'
$ws.Range("E33").Value = '  -6.14%  '
$ws.Range("D34").Value = '0.08881'
$ws.Range("E34").Value = '  -4.14%  '
$ws.Range("D35").Value = '12.20'
$ws.Range("E35").Value = '  -6.86%  '
$ws.Range("D36").Value = '0.02276'
$ws.Range("E36").Value = '  -3.41%  '
$ws.Range("D37").Value = '0.2095'
$ws.Range("E37").Value = '  -4.13%  '
$ws.Range("D38").Value = '0.06002'
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("D39").Value = '0.6336'
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("D40").Value = '4.946'
$ws.Range("E40").Value = '  -4.56%  '
$ws.Range("D41").Value = '1.182'
$ws.Range("E41").Value = '  -1.32%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").Value = '7.877'
$ws.Range("E43").Value = '  -3.35%  '
$ws.Range("D44").Value = '1.385'
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("D45").Value = '13.33'
$ws.Range("E45").Value = '  -4.22%  '
$ws.Range("D46").Value = '0.5857'
$ws.Range("E46").Value = '  -3.82%  '
$ws.Range("D47").Value = '3.685'
$ws.Range("D48").Value = '1.972'
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").Value = '122.15'
$ws.Range("E49").Value = '  -3.81%  '
$ws.Range("D50").Value = '1.171'
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("D51").Value = '0.06805'
